# Rows 6-12 of the "Artfynd" sheet get their record data rotated in a
# single 7-cycle: 6<-11<-12<-10<-9<-8<-7<-6 (i.e. new row6 = old row11,
# new row7 = old row6, new row8 = old row7, new row9 = old row8,
# new row10 = old row9, new row11 = old row12, new row12 = old row10).
# Columns C, I, T, U, V, W, Z, AB, AD, AE, AG, AT are identical across all
# these rows already, so only the cells below actually change value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6  (<= old row 11)
$ws.Cells.Item(6,1).Value  = 74825621            # A6  Id
$ws.Cells.Item(6,17).Value = 635778.1134565037   # Q6  Ost
$ws.Cells.Item(6,18).Value = 6517624.051674079   # R6  Nord

# Row 7  (<= old row 6)
$ws.Cells.Item(7,1).Value  = 74825618                      # A7  Id
$ws.Cells.Item(7,2).Value  = 4717                           # B7  Taxonsorteringsordning
$ws.Cells.Item(7,5).Value  = 102306                         # E7  TaxonId
$ws.Cells.Item(7,6).Value  = "Granbarkgnagare"              # F7  Artnamn
$ws.Cells.Item(7,7).Value  = "Microbregma emarginatum"      # G7  Vetenskapligt namn
$ws.Cells.Item(7,8).Value  = "(Duftschmid, 1825)"           # H7  Auktor
$ws.Cells.Item(7,17).Value = 635969.0379152392              # Q7  Ost
$ws.Cells.Item(7,18).Value = 6517709.110210423              # R7  Nord

# Row 8  (<= old row 7)
$ws.Cells.Item(8,1).Value  = 74825623                  # A8  Id
$ws.Cells.Item(8,2).Value  = 93158                      # B8  Taxonsorteringsordning
$ws.Cells.Item(8,4).Value  = "LC"                       # D8  Rödlistade
$ws.Cells.Item(8,5).Value  = 2818                       # E8  TaxonId
$ws.Cells.Item(8,6).Value  = "Stubbspretmossa"          # F8  Artnamn
$ws.Cells.Item(8,7).Value  = "Herzogiella seligeri"     # G8  Vetenskapligt namn
$ws.Cells.Item(8,8).Value  = "(Brid.) Z.Iwats."         # H8  Auktor
$ws.Cells.Item(8,17).Value = 635877.0551628179          # Q8  Ost
$ws.Cells.Item(8,18).Value = 6517560.061384213          # R8  Nord

# Row 9  (<= old row 8)
$ws.Cells.Item(9,1).Value  = 74825619                          # A9  Id
$ws.Cells.Item(9,2).Value  = 89410                              # B9  Taxonsorteringsordning
$ws.Cells.Item(9,5).Value  = 5432                               # E9  TaxonId
$ws.Cells.Item(9,6).Value  = "Granticka"                        # F9  Artnamn
$ws.Cells.Item(9,7).Value  = "Porodaedalea chrysoloma"          # G9  Vetenskapligt namn
$ws.Cells.Item(9,8).Value  = "(Fr.) Fiasson & Niemelä"          # H9  Auktor
$ws.Cells.Item(9,16).Value = "Hånö-Kärr, Srm"                   # P9  Lokalnamn
$ws.Cells.Item(9,17).Value = 635969.0379152392                  # Q9  Ost
$ws.Cells.Item(9,18).Value = 6517709.110210423                  # R9  Nord
$ws.Cells.Item(9,19).Value = 10                                 # S9  Noggrannhet
# Y9/AA9 are text dates ("yyyy-mm-dd"), not real date serials - force
# text format first so Excel doesn't silently convert them to dates.
$ws.Cells.Item(9,25).NumberFormat = "@"
$ws.Cells.Item(9,25).Value = "2018-09-12"                       # Y9  Startdatum
$ws.Cells.Item(9,27).NumberFormat = "@"
$ws.Cells.Item(9,27).Value = "2018-09-12"                       # AA9 Slutdatum
$ws.Cells.Item(9,49).Value = "Markus Forsberg"                  # AW9 Rapportör
$ws.Cells.Item(9,50).Value = "Markus Forsberg"                  # AX9 Observatörer
$ws.Cells.Item(9,51).Value = ""                                 # AY9 Projektnamn (cleared)

# Row 10  (<= old row 9)
$ws.Cells.Item(10,1).Value  = 84982788              # A10  Id
$ws.Cells.Item(10,2).Value  = 90138                  # B10  Taxonsorteringsordning
$ws.Cells.Item(10,4).Value  = "NT"                   # D10  Rödlistade
$ws.Cells.Item(10,5).Value  = 366                    # E10  TaxonId
$ws.Cells.Item(10,6).Value  = "Kandelabersvamp"      # F10  Artnamn
$ws.Cells.Item(10,7).Value  = "Artomyces pyxidatus"  # G10  Vetenskapligt namn
$ws.Cells.Item(10,8).Value  = "(Pers.) Jülich"       # H10  Auktor
$ws.Cells.Item(10,17).Value = 635883.4449681807      # Q10  Ost
$ws.Cells.Item(10,18).Value = 6517466.875374788      # R10  Nord

# Row 11  (<= old row 12)
$ws.Cells.Item(11,1).Value  = 84982787                           # A11  Id
$ws.Cells.Item(11,2).Value  = 43464                               # B11  Taxonsorteringsordning
$ws.Cells.Item(11,5).Value  = 101735                              # E11  TaxonId
$ws.Cells.Item(11,6).Value  = "Jättesvampmal"                    # F11  Artnamn
$ws.Cells.Item(11,7).Value  = "Scardia boletella"                # G11  Vetenskapligt namn
$ws.Cells.Item(11,8).Value  = "(Fabricius, 1794)"                # H11  Auktor
$ws.Cells.Item(11,16).Value = "Kärr, Srm"                        # P11  Lokalnamn
$ws.Cells.Item(11,17).Value = 635778.11899246                    # Q11  Ost
$ws.Cells.Item(11,18).Value = 6517579.423625848                  # R11  Nord
$ws.Cells.Item(11,19).Value = 5                                  # S11  Noggrannhet
$ws.Cells.Item(11,25).NumberFormat = "@"
$ws.Cells.Item(11,25).Value = "2019-11-15"                       # Y11  Startdatum
$ws.Cells.Item(11,27).NumberFormat = "@"
$ws.Cells.Item(11,27).Value = "2019-11-15"                       # AA11 Slutdatum
$ws.Cells.Item(11,49).Value = "Ralf Lundmark"                    # AW11 Rapportör
$ws.Cells.Item(11,50).Value = "Ralf Lundmark"                    # AX11 Observatörer
$ws.Cells.Item(11,51).Value = "Kryptogamer i Södermanlands län"  # AY11 Projektnamn

# Row 12  (<= old row 10)
$ws.Cells.Item(12,1).Value  = 84982789                 # A12  Id
$ws.Cells.Item(12,2).Value  = 93158                     # B12  Taxonsorteringsordning
$ws.Cells.Item(12,5).Value  = 2818                      # E12  TaxonId
$ws.Cells.Item(12,6).Value  = "Stubbspretmossa"         # F12  Artnamn
$ws.Cells.Item(12,7).Value  = "Herzogiella seligeri"    # G12  Vetenskapligt namn
$ws.Cells.Item(12,8).Value  = "(Brid.) Z.Iwats."        # H12  Auktor
$ws.Cells.Item(12,17).Value = 635878.212690325          # Q12  Ost
$ws.Cells.Item(12,18).Value = 6517571.518665905         # R12  Nord
